$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.558.49"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.576.44"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'288.64"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").Value = "'0.3684"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").Value = "'48.53"
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "'0.07473"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").Value = "'6.009"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'6.953"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "1.582.33"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "'0.00001118"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "'88.72"
$ws.Range("D19").Value = "'0.06755"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'6.432"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'16.58"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "'12.19"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "22.560.00"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "'2.407"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").Value = "'2.614"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'152.73"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").Value = "'19.67"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "'5.023"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "'124.43"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "1.757.43"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "'1.071"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'6.207"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'2.002"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "'9.759"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").Value = "'0.08294"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "'0.02460"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "'0.2276"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").Value = "'5.456"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'0.06399"
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "'11.40"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'0.6356"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'14.00"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").Value = "'0.6186"
$ws.Range("E46").Value = "  +5.45%  "
$ws.Range("D47").Value = "'3.770"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "'2.062"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").Value = "'125.36"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").Value = "'1.221"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'0.07281"
$ws.Range("E51").Value = "  -0.78%  "
